# Rebuild the LR-pairs data table with the re-run NATMI results
# (adds the "ECs" sending-cluster group and refreshes all computed values,
# per "Natmi following Dr Hou advice").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Sema4g", "Plxnb2", "ECs", 2, 0.6666666666666666, 0.172848, 0.518544, 0.04011888417844697, 0.04011888417844697, 3, 1, 8.777551666666666, 26.332655, 0.07454818073713242, 0.07454818073713242, 1.51718225048, 13.65464025432, 0.002990789828706947, 0.002990789828706947),
    @("ECs", "Sema4g", "Plxnb2", "FAPs", 2, 0.6666666666666666, 0.172848, 0.518544, 0.04011888417844697, 0.04011888417844697, 3, 1, 36.14140700000001, 108.424221, 0.306950758417288, 0.306950758417288, 6.246969917136002, 56.22272925422401, 0.01231452192542963, 0.01231452192542963),
    @("ECs", "Sema4g", "Plxnb2", "M2", 2, 0.6666666666666666, 0.172848, 0.518544, 0.04011888417844697, 0.04011888417844697, 3, 1, 54.580447, 163.741341, 0.4635544377507104, 0.4635544377507104, 9.434121103056, 84.907089927504, 0.01859728679852585, 0.01859728679852585),
    @("ECs", "Sema4g", "Plxnb2", "sCs", 2, 0.6666666666666666, 0.172848, 0.518544, 0.04011888417844697, 0.04011888417844697, 3, 1, 18.24393266666667, 54.731798, 0.1549466230948692, 0.1549466230948692, 3.153427273568, 28.380845462112, 0.006216285625784535, 0.006216285625784533),
    @("FAPs", "Sema4g", "Plxnb2", "ECs", 3, 1, 2.436670666666667, 7.310012, 0.5655634329411919, 0.5655634329411919, 3, 1, 8.777551666666666, 26.332655, 0.07454818073713242, 0.07454818073713242, 21.38800267131777, 192.49202404186, 0.04216172501721305, 0.04216172501721305),
    @("FAPs", "Sema4g", "Plxnb2", "FAPs", 3, 1, 2.436670666666667, 7.310012, 0.5655634329411919, 0.5655634329411919, 3, 1, 36.14140700000001, 108.424221, 0.306950758417288, 0.306950758417288, 88.06470628896135, 792.582356600652, 0.1736001246743839, 0.1736001246743838),
    @("FAPs", "Sema4g", "Plxnb2", "M2", 3, 1, 2.436670666666667, 7.310012, 0.5655634329411919, 0.5655634329411919, 3, 1, 54.580447, 163.741341, 0.4635544377507104, 0.4635544377507104, 132.9945741784547, 1196.951167606092, 0.2621694391694158, 0.2621694391694158),
    @("FAPs", "Sema4g", "Plxnb2", "sCs", 3, 1, 2.436670666666667, 7.310012, 0.5655634329411919, 0.5655634329411919, 3, 1, 18.24393266666667, 54.731798, 0.1549466230948692, 0.1549466230948692, 44.45445557350845, 400.090100161576, 0.08763214408017922, 0.0876321440801792),
    @("sCs", "Sema4g", "Plxnb2", "ECs", 3, 1, 1.698876333333333, 5.096629, 0.3943176828803611, 0.3943176828803611, 3, 1, 8.777551666666666, 26.332655, 0.07454818073713242, 0.07454818073713242, 14.91197479111056, 134.207773119995, 0.02939566589121242, 0.02939566589121242),
    @("sCs", "Sema4g", "Plxnb2", "FAPs", 3, 1, 1.698876333333333, 5.096629, 0.3943176828803611, 0.3943176828803611, 3, 1, 36.14140700000001, 108.424221, 0.306950758417288, 0.306950758417288, 61.39978100566768, 552.5980290510091, 0.1210361118174745, 0.1210361118174745),
    @("sCs", "Sema4g", "Plxnb2", "M2", 3, 1, 1.698876333333333, 5.096629, 0.3943176828803611, 0.3943176828803611, 3, 1, 54.580447, 163.741341, 0.4635544377507104, 0.4635544377507104, 92.72542967105434, 834.528867039489, 0.1827877117827687, 0.1827877117827687),
    @("sCs", "Sema4g", "Plxnb2", "sCs", 3, 1, 1.698876333333333, 5.096629, 0.3943176828803611, 0.3943176828803611, 3, 1, 18.24393266666667, 54.731798, 0.1549466230948692, 0.1549466230948692, 30.99418543432689, 278.947668908942, 0.06109819338890549, 0.06109819338890547)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
